$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "npm " + "init" (two runs, wrapped by <w:proofErr spellStart/spellEnd>)
#    -> single run "npm init" with no proofErr markers.
#    A plain Find/Replace merges the two runs and drops the spellStart
#    marker, but leaves an orphaned spellEnd behind. To get a fully clean
#    paragraph we rebuild it: insert a fresh paragraph (which inherits the
#    same paragraph/run formatting), give it the finished text, then delete
#    the old (proofErr-laden) paragraph entirely.
# ---------------------------------------------------------------------------
$oldPara = $d.Paragraphs.Item(2)
$insPoint = $d.Range($oldPara.Range.Start, $oldPara.Range.Start)
$insPoint.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item(2)
$newPara.Range.Text = "npm init"

$oldPara = $d.Paragraphs.Item(3)
$oldRange = $d.Range($oldPara.Range.Start, $oldPara.Range.End)
$oldRange.Delete()

# ---------------------------------------------------------------------------
# 2) & 4) Mark the runs that hold the two inline images as NoProofing, which
#    emits <w:rPr><w:noProof/></w:rPr> on their runs.
# ---------------------------------------------------------------------------
$shapes = $d.InlineShapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shapes.Item($i).Range.NoProofing = $true
}

# ---------------------------------------------------------------------------
# 3) Drop the stale <w:lastRenderedPageBreak/> cached before "You are now
#    able to install the rest" - a no-op Find/Replace over the run's text
#    clears the render cache marker while leaving the text/run untouched.
# ---------------------------------------------------------------------------
$target = "You are now able to install the rest"
$d.Content.Find.Execute($target, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $target, 2)
